$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A updates (the "popularity"/unused-index column)
$ws.Range("A2").Value = 27
$ws.Range("A3").Value = 16
$ws.Range("A4").Value = 37
$ws.Range("A5").Value = 50
$ws.Range("A9").Value = 38
$ws.Range("A11").Value = 41
$ws.Range("A13").Value = 26
$ws.Range("A18").Value = 58
$ws.Range("A19").Value = 42
$ws.Range("A20").Value = 8
$ws.Range("A26").Value = 10
$ws.Range("A27").Value = 9
$ws.Range("A29").Value = 57
$ws.Range("A30").Value = 15
$ws.Range("A35").Value = 30
$ws.Range("A37").Value = 52
$ws.Range("A39").Value = 43
$ws.Range("A40").Value = 40
$ws.Range("A42").Value = 28
$ws.Range("A43").Value = 39
$ws.Range("A44").Value = 33
$ws.Range("A45").Value = 51
$ws.Range("A49").Value = 17
$ws.Range("A53").Value = 4
$ws.Range("A54").Value = 11
$ws.Range("A55").Value = 31
$ws.Range("A56").Value = 32
$ws.Range("A58").Value = 6

# Rating/user_ratings_total updates
$ws.Range("E9").Value = 473
$ws.Range("E22").Value = 528
$ws.Range("D45").Value = 5
$ws.Range("E45").Value = 62
$ws.Range("E49").Value = 845
